$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold plain-text values (coin prices / % changes / names /
# links) stored as inline strings in the original workbook, e.g. "1.006"
# or "11.00". A bare numeric-looking string assigned via .Value gets
# auto-coerced to a real number by Excel/COM (dropping trailing zeros,
# switching to scientific notation, etc.), which would not match the
# source data. Prefixing with a leading apostrophe forces Excel to treat
# the assignment as literal text (exactly like typing `'1.006` into a
# cell in the UI) while leaving the cell format untouched ("General").

$ws.Range("D2").Value = "'26.597.10"
$ws.Range("E2").Value = "'  -7.23%  "

$ws.Range("D3").Value = "'1.694.09"
$ws.Range("E3").Value = "'  -5.99%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "'  +0.33%  "

$ws.Range("D5").Value = "'219.95"
$ws.Range("E5").Value = "'  -5.13%  "

$ws.Range("D6").Value = "'0.5114"
$ws.Range("E6").Value = "'  -13.79%  "

$ws.Range("E7").Value = "'  +0.20%  "

$ws.Range("D8").Value = "'0.2655"
$ws.Range("E8").Value = "'  -4.33%  "

$ws.Range("D9").Value = "'22.21"
$ws.Range("E9").Value = "'  -4.72%  "

$ws.Range("D10").Value = "'0.06291"
$ws.Range("E10").Value = "'  -7.71%  "

$ws.Range("D11").Value = "'0.07358"
$ws.Range("E11").Value = "'  -1.99%  "

$ws.Range("D12").Value = "'1.702.31"

$ws.Range("D13").Value = "'4.514"
$ws.Range("E13").Value = "'  -5.30%  "

$ws.Range("D14").Value = "'0.5852"
$ws.Range("E14").Value = "'  -5.96%  "

$ws.Range("D15").Value = "'1.925.96"
$ws.Range("E15").Value = "'  -5.90%  "

$ws.Range("D16").Value = "'0.000008405"
$ws.Range("E16").Value = "'  -8.46%  "

$ws.Range("D17").Value = "'65.59"
$ws.Range("E17").Value = "'  -13.20%  "

$ws.Range("D18").Value = "'26.629.88"
$ws.Range("E18").Value = "'  -7.06%  "

$ws.Range("D19").Value = "'5.021"
$ws.Range("E19").Value = "'  -8.15%  "

$ws.Range("E20").Value = "'  +0.19%  "

$ws.Range("D21").Value = "'11.00"
$ws.Range("E21").Value = "'  -4.36%  "

$ws.Range("D22").Value = "'186.39"
$ws.Range("E22").Value = "'  -11.60%  "

$ws.Range("D23").Value = "'6.279"
$ws.Range("E23").Value = "'  -8.06%  "

$ws.Range("E24").Value = "'  +0.22%  "

$ws.Range("D25").Value = "'144.65"
$ws.Range("E25").Value = "'  -5.83%  "

$ws.Range("D26").Value = "'7.527"
$ws.Range("E26").Value = "'  -4.13%  "

$ws.Range("E27").Value = "'  -8.81%  "

$ws.Range("D28").Value = "'15.61"
$ws.Range("E28").Value = "'  -5.07%  "

$ws.Range("D29").Value = "'1.334"
$ws.Range("E29").Value = "'  -5.17%  "

$ws.Range("D30").Value = "'0.05671"
$ws.Range("E30").Value = "'  -7.83%  "

$ws.Range("D31").Value = "'1.338"
$ws.Range("E31").Value = "'  -6.12%  "

$ws.Range("D32").Value = "'3.519"
$ws.Range("E32").Value = "'  -6.79%  "

$ws.Range("D33").Value = "'3.497"
$ws.Range("E33").Value = "'  -6.57%  "

$ws.Range("D34").Value = "'1.659"
$ws.Range("E34").Value = "'  -4.00%  "

$ws.Range("D35").Value = "'1.021"
$ws.Range("E35").Value = "'  -3.45%  "

$ws.Range("D36").Value = "'0.6051"
$ws.Range("E36").Value = "'  -5.70%  "

$ws.Range("D37").Value = "'2.362"
$ws.Range("E37").Value = "'  -5.57%  "

$ws.Range("D38").Value = "'2.680"
$ws.Range("E38").Value = "'  -1.35%  "

$ws.Range("B39").Value = "'VeChain"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01612"
$ws.Range("E39").Value = "'  -4.69%  "

$ws.Range("B40").Value = "'Maker"
$ws.Range("C40").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'1.099.73"
$ws.Range("E40").Value = "'  -4.22%  "

$ws.Range("D41").Value = "'0.8663"
$ws.Range("E41").Value = "'  -2.36%  "

$ws.Range("D42").Value = "'5.866"
$ws.Range("E42").Value = "'  -10.45%  "

$ws.Range("E43").Value = "'  -0.29%  "

$ws.Range("D44").Value = "'99.17"
$ws.Range("E44").Value = "'  -0.91%  "

$ws.Range("D45").Value = "'1.858.33"
$ws.Range("E45").Value = "'  -4.98%  "

$ws.Range("D46").Value = "'0.00000000110"
$ws.Range("E46").Value = "'  -2.80%  "

$ws.Range("D47").Value = "'56.79"

$ws.Range("D48").Value = "'8.155"
$ws.Range("E48").Value = "'  -2.55%  "

$ws.Range("E49").Value = "'  +0.40%  "

$ws.Range("D50").Value = "'0.05251"
$ws.Range("E50").Value = "'  -4.01%  "

$ws.Range("D51").Value = "'0.4331"
$ws.Range("E51").Value = "'  -3.21%  "
